$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values, updated per the automatic electricity price refresh.
$ws.Range("A2").Value = 45906
$ws.Range("B2").Value = 74.64
$ws.Range("C2").Value = 58.78
$ws.Range("D2").Value = 48.71
$ws.Range("E2").Value = 39.9
$ws.Range("F2").Value = 34.8
$ws.Range("G2").Value = 33.46
$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 42
$ws.Range("J2").Value = 34.9
$ws.Range("K2").Value = 15.23
$ws.Range("L2").Value = 1.72
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = -0.02
$ws.Range("O2").Value = -0.86
$ws.Range("P2").Value = -0.99
$ws.Range("Q2").Value = -0.16
$ws.Range("R2").Value = -0.02
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 24.37
$ws.Range("V2").Value = 42
$ws.Range("W2").Value = 58.54
$ws.Range("X2").Value = 55.44
$ws.Range("Y2").Value = 39.5
$ws.Range("Z2").Value = 26.54
$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 55.51
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 66.70999999999999
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 50.27
$ws.Range("AG2").Value = "9h-19h"
